# 2363-RBI-EPP... Newcreateloan1.xlsx -- "modified test cases on overdue fix"
#
# Summary of the change (per the commit's OOXML diff):
#  - Summary sheet ("Summary"): shrink used range from G10 to G5 (rows 7-10 were
#    scratch/blank rows, now removed), tweak a couple of rounded values, and
#    drop the custom "0.00" cell style (style index 21) in favour of the plain
#    style already used by neighbouring cells.
#  - "Repayment schedule": the same style-21 cleanup, several amortisation
#    values nudged by a cent (rounding fix), and a 7th instalment row added
#    (row 9) because the loan was not fully settled by instalment 6 anymore.
#  - Sheet selections / active cell bookkeeping updated to match, and the
#    active tab moves from "Repayment schedule" to "Transactions".

$wb = $excel.ActiveWorkbook

$wsInput   = $wb.Worksheets.Item("NewLoanInput")
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSched   = $wb.Worksheets.Item("Repayment schedule")
$wsTxn     = $wb.Worksheets.Item("Transactions")

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------

# Re-style A3/E3/F3 away from the soon-to-be-orphaned "0.00" style (21) onto
# the plain style already used by B3 (style 14) -- format only, values are
# rewritten right after.
$wsSummary.Range("B3").Copy() | Out-Null
$wsSummary.Range("A3").PasteSpecial(-4122) | Out-Null
$wsSummary.Range("B3").Copy() | Out-Null
$wsSummary.Range("E3").PasteSpecial(-4122) | Out-Null
$wsSummary.Range("B3").Copy() | Out-Null
$wsSummary.Range("F3").PasteSpecial(-4122) | Out-Null

$wsSummary.Range("F2").Value2 = 1666.66
$wsSummary.Range("A3").Value2 = 74.239999999999995
$wsSummary.Range("E3").Value2 = 74.239999999999995
$wsSummary.Range("F3").Value2 = 28.08

# Drop the now-unused scratch rows 7:10.
$wsSummary.Rows("7:10").Delete() | Out-Null

# ---------------------------------------------------------------------------
# Repayment schedule
# ---------------------------------------------------------------------------

# Re-style the data block (rows 2-8) off the "0.00" style (21) onto the plain
# style already used elsewhere on the row (style 14), then re-apply the
# special columns (C = date, G = numeric amount columns) on top.
$wsSched.Range("I2").Copy() | Out-Null
$wsSched.Range("A2:P8").PasteSpecial(-4122) | Out-Null

$wsSched.Range("C2").Copy() | Out-Null
$wsSched.Range("C2:C8").PasteSpecial(-4122) | Out-Null

$wsSummary.Range("A2").Copy() | Out-Null
$wsSched.Range("G2").PasteSpecial(-4122) | Out-Null

$wsSummary.Range("F2").Copy() | Out-Null
$wsSched.Range("G3:G6").PasteSpecial(-4122) | Out-Null

# Rounding-cent fixes on the existing instalment rows.
$wsSched.Range("G4").Value2 = 3333.34
$wsSched.Range("K4").Value2 = 848.67
$wsSched.Range("P4").Value2 = 848.67

$wsSched.Range("G5").Value2 = 2500.0100000000002
$wsSched.Range("K5").Value2 = 854.56
$wsSched.Range("P5").Value2 = 854.56

$wsSched.Range("G6").Value2 = 1666.68

$wsSched.Range("G7").Value2 = 833.35
$wsSched.Range("K7").Value2 = 841.82
$wsSched.Range("P7").Value2 = 841.82

$wsSched.Range("G8").Value2 = 0.02

# New instalment row (#7), fully formatted/styled off row 8.
$wsSched.Range("A8:P8").Copy() | Out-Null
$wsSched.Range("A9:P9").PasteSpecial(-4122) | Out-Null

$wsSched.Range("A9").Value2 = 7
$wsSched.Range("B9").Value2 = 31
$wsSched.Range("C9").Value2 = 42217
$wsSched.Range("D9").Value2 = ""
$wsSched.Range("E9").Value2 = ""
$wsSched.Range("F9").Value2 = 0.02
$wsSched.Range("G9").Value2 = 0
$wsSched.Range("H9").Value2 = 0
$wsSched.Range("I9").Value2 = 0
$wsSched.Range("J9").Value2 = 0
$wsSched.Range("K9").Value2 = 0.02
$wsSched.Range("L9").Value2 = 0
$wsSched.Range("M9").Value2 = 0
$wsSched.Range("N9").Value2 = 0
$wsSched.Range("O9").Value2 = 0
$wsSched.Range("P9").Value2 = 0.02

# ---------------------------------------------------------------------------
# Selections / active tab bookkeeping
# ---------------------------------------------------------------------------

$wsSummary.Range("A7:XFD15").Select() | Out-Null
$wsSched.Range("A10:XFD10").Select() | Out-Null
$wsTxn.Activate() | Out-Null
$wsTxn.Range("H13").Select() | Out-Null
